$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Force these cells to keep storing plain text (matches the source file,
# where every value in this sheet is an inline/shared string, not a
# number or date) instead of letting Excel auto-convert on assignment.
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B13").NumberFormat = "@"

# Total Forecast (16 Weeks): 13 -> 17
$ws.Range("B9").Value = "17"

# Total Forecast (8 Weeks): 7 -> 9
$ws.Range("B10").Value = "9"

# Total Forecast (4 Weeks): 3 -> 4
$ws.Range("B11").Value = "4"

# Max Forecast Week: 2024-12-29 -> 2024-12-15
$ws.Range("B13").Value = "2024-12-15"
